$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (data that ends up in destination row
# was previously found in source row), for columns D, L, M, N, O, P, Q, R, S, T.
$rowMap = @{
    2  = 11
    3  = 8
    4  = 6
    5  = 10
    6  = 2
    7  = 13
    8  = 14
    9  = 3
    10 = 16
    11 = 15
    12 = 9
    13 = 12
    14 = 4
    15 = 7
    16 = 5
}

# Capture the "before" values for the columns that move, keyed by row number.
# NOTE: use .Value2 for reads - .Value getter is unreliable in this runtime.
$before = @{}
foreach ($r in 2..16) {
    $before[$r] = @{
        D = $ws.Range("D$r").Value2
        L = $ws.Range("L$r").Value2
        M = $ws.Range("M$r").Value2
        N = $ws.Range("N$r").Value2
        O = $ws.Range("O$r").Value2
        P = $ws.Range("P$r").Value2
        Q = $ws.Range("Q$r").Value2
        R = $ws.Range("R$r").Value2
        S = $ws.Range("S$r").Value2
        T = $ws.Range("T$r").Value2
    }
}

# Apply the permutation: for each destination row, pull the values from the
# recorded source row's original data.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $src = $before[$srcRow]

    $ws.Range("D$destRow").Value = $src.D
    $ws.Range("L$destRow").Value = $src.L
    $ws.Range("M$destRow").Value = $src.M
    $ws.Range("N$destRow").Value = $src.N
    $ws.Range("O$destRow").Value = $src.O
    $ws.Range("P$destRow").Value = $src.P
    $ws.Range("Q$destRow").Value = $src.Q
    $ws.Range("R$destRow").Value = $src.R
    $ws.Range("S$destRow").Value = $src.S
    $ws.Range("T$destRow").Value = $src.T
}
